$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.347.62'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.10%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.906.88'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.10%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '487.62'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.45%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.73'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.59%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.622'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.40%  '
$ws.Range("E9").Value = '  +2.85%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.180'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +8.78%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0000357'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.60%  '
$ws.Range("E12").Value = '  +1.44%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '10.50'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.47%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.520.04'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.36%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.889.60'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.05%  '
$ws.Range("E16").Value = '  -2.56%  '
$ws.Range("E17").Value = '  -0.46%  '
$ws.Range("E18").Value = '  +1.65%  '
$ws.Range("E19").Value = '  +1.92%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '68.451.61'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.09%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '433.57'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.31%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.56'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +5.77%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '14.70'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.16%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '90.05'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.63%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.36'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +18.23%  '
$ws.Range("E26").Value = '  +4.79%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.99'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -6.07%  '
$ws.Range("E28").Value = '  -1.39%  '
$ws.Range("E29").Value = '  -2.16%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '712.40'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.67%  '
$ws.Range("E31").Value = '  +0.71%  '
$ws.Range("E32").Value = '  +0.62%  '
$ws.Range("E33").Value = '  +2.21%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0₃0899'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.84%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.08'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +6.63%  '
$ws.Range("E36").Value = '  +3.74%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '40.83'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.77%  '
$ws.Range("E38").Value = '  +19.44%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.148'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.63%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.00'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.16%  '
$ws.Range("E41").Value = '  +4.63%  '
$ws.Range("E42").Value = '  +7.64%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.10'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.18%  '
$ws.Range("E44").Value = '  -1.25%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0₆0377'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +28.45%  '
$ws.Range("E46").Value = '  +1.24%  '
$ws.Range("E47").Value = '  +8.01%  '
$ws.Range("E48").Value = '  -0.01%  '
$ws.Range("E49").Value = '  -1.44%  '
$ws.Range("E50").Value = '  -1.79%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '142.41'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.22%  '
